# Add 12 new login rows (3 mobile numbers sharing one identity, x4 groups)
# to Sheet1, and mark the corresponding mobile numbers as "used" on the
# "Test Data" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)        # "Sheet1"
$ws2 = $wb.Worksheets.Item(2)        # "Test Data"

# --- New rows for Sheet1 (row -> Mobile, UserName, Email, Password) ---
$newRows = @(
    @("7980000069", "test68629",  "test68629@gmail.com",  "SoftSuave6273"),
    @("7980000070", "test68629",  "test68629@gmail.com",  "SoftSuave6273"),
    @("7980000072", "test68629",  "test68629@gmail.com",  "SoftSuave6273"),
    @("7980000073", "test44642",  "test44642@gmail.com",  "SoftSuave101694"),
    @("7980000074", "test44642",  "test44642@gmail.com",  "SoftSuave101694"),
    @("7980000075", "test44642",  "test44642@gmail.com",  "SoftSuave101694"),
    @("7980000076", "test86530",  "test86530@gmail.com",  "SoftSuave105706"),
    @("7980000077", "test86530",  "test86530@gmail.com",  "SoftSuave105706"),
    @("7980000078", "test86530",  "test86530@gmail.com",  "SoftSuave105706"),
    @("7980000079", "test14552",  "test14552@gmail.com",  "SoftSuave148393"),
    @("7980000080", "test14552",  "test14552@gmail.com",  "SoftSuave148393"),
    @("7980000081", "test14552",  "test14552@gmail.com",  "SoftSuave148393")
)

$startRow = 41
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $mobile = $newRows[$i][0]

    # The mobile number cell has to stay text (shared string), like every
    # other mobile-number cell already on this sheet. A plain .Value=
    # assignment of a numeric-looking string is auto-coerced to a number,
    # so instead write a literal-text formula and freeze it back down to a
    # plain value via copy / paste-special (values only).
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Formula = '="' + $mobile + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    $ws1.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws1.Cells.Item($r, 3).Value = $newRows[$i][2]
    $ws1.Cells.Item($r, 4).Value = $newRows[$i][3]
}
$excel.CutCopyMode = $false

# --- Mark mobile numbers 7980000069..7980000081 as "used" on Test Data ---
for ($r = 70; $r -le 82; $r++) {
    $ws2.Cells.Item($r, 2).Value = "used"
}
